$d = $word.ActiveDocument

# The title paragraph gains a leading space.
$titleStart = $d.Range(0, 0)
$titleStart.InsertBefore(" ")

# The list item paragraph loses its leading space and its trailing
# colon becomes a period.
$find = " Ներմուծում ենք բազմանկյան n հատ գագաթների կոորդինատները և պահում զանգվածում  XandY[i][j]  i=1..2, j=1...n:"
$replace = "Ներմուծում ենք բազմանկյան n հատ գագաթների կոորդինատները և պահում զանգվածում  XandY[i][j]  i=1..2, j=1...n."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
